# Updates cryptos list price (D) and 1h volume change (E) columns to the
# latest scrape for each coin row, per the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as literal text (e.g. "59.654.98", "558.94") so the
# odd multi-dot big-number formatting round-trips unchanged. Plain decimal
# strings like "558.94" would otherwise be auto-converted to a Double by the
# Value setter, so prefix with the Excel text-prefix apostrophe and then put
# the cell style back to Normal (the apostrophe alone leaves quotePrefix set).
$apostrophe = "'"

$ws.Range("D2").Value = $apostrophe + "59.654.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = $apostrophe + "2.367.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.63%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = $apostrophe + "558.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").Value = $apostrophe + "137.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = $apostrophe + "0.527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = $apostrophe + "2.363.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").Value = $apostrophe + "0.105"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = $apostrophe + "25.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").Value = $apostrophe + "2.794.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("E16").Value = "  -3.18%  "
$ws.Range("D17").Value = $apostrophe + "59.721.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").Value = $apostrophe + "2.365.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = $apostrophe + "8.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.41%  "
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").Value = $apostrophe + "321.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = $apostrophe + "5.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -3.69%  "
$ws.Range("D26").Value = $apostrophe + "64.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("D27").Value = $apostrophe + "558.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.77%  "
$ws.Range("D28").Value = $apostrophe + "8.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.45%  "
$ws.Range("D29").Value = $apostrophe + "2.483.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.16%  "
$ws.Range("D30").Value = $apostrophe + "0.0₃0919"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").Value = $apostrophe + "7.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = $apostrophe + "1.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.04%  "
$ws.Range("D37").Value = $apostrophe + "152.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.89%  "
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").Value = $apostrophe + "4.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").Value = $apostrophe + "18.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").Value = $apostrophe + "4.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").Value = $apostrophe + "2.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("D46").Value = $apostrophe + "0.0₆0296"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.26%  "
$ws.Range("D47").Value = $apostrophe + "138.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = $apostrophe + "0.584"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = $apostrophe + "19.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.52%  "
